$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) values are stored as text, matching the source data
# (values like "1.002" or "321.06" would otherwise be auto-converted to numbers)

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '29.228.71'
$ws.Cells.Item(2, 5).Value = '  +0.97%  '

$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '1.909.71'
$ws.Cells.Item(3, 5).Value = '  +1.34%  '

$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '1.002'
$ws.Cells.Item(4, 5).Value = '  +0.10%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '321.06'
$ws.Cells.Item(5, 5).Value = '  -2.91%  '

$ws.Cells.Item(6, 5).Value = '  +0.04%  '

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.4712'
$ws.Cells.Item(7, 5).Value = '  +2.77%  '

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.4063'
$ws.Cells.Item(8, 5).Value = '  -0.06%  '

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '47.75'
$ws.Cells.Item(9, 5).Value = '  +0.47%  '

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.08034'
$ws.Cells.Item(10, 5).Value = '  +0.92%  '

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.9997'
$ws.Cells.Item(11, 5).Value = '  +1.00%  '

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '22.44'
$ws.Cells.Item(12, 5).Value = '  +3.73%  '

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '1.920.61'
$ws.Cells.Item(13, 5).Value = '  +2.40%  '

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '5.877'
$ws.Cells.Item(14, 5).Value = '  -0.38%  '

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '7.117'
$ws.Cells.Item(15, 5).Value = '  +1.05%  '

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '89.60'
$ws.Cells.Item(16, 5).Value = '  +1.46%  '

$ws.Cells.Item(17, 5).Value = '  +0.04%  '

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '0.06634'
$ws.Cells.Item(18, 5).Value = '  +1.26%  '

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '0.00001029'
$ws.Cells.Item(19, 5).Value = '  +0.23%  '

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '17.66'
$ws.Cells.Item(20, 5).Value = '  +1.51%  '

$ws.Cells.Item(21, 5).Value = '  +0.09%  '

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '29.260.18'
$ws.Cells.Item(22, 5).Value = '  +0.87%  '

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '5.512'
$ws.Cells.Item(23, 5).Value = '  +1.88%  '

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '11.43'
$ws.Cells.Item(24, 5).Value = '  +1.66%  '

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '2.197'
$ws.Cells.Item(25, 5).Value = '  -0.47%  '

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '2.175.59'
$ws.Cells.Item(26, 5).Value = '  +3.49%  '

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '155.06'
$ws.Cells.Item(27, 5).Value = '  -1.06%  '

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '19.80'
$ws.Cells.Item(28, 5).Value = '  +1.28%  '

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '5.992'
$ws.Cells.Item(29, 5).Value = '  +11.06%  '

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '2.102'
$ws.Cells.Item(30, 5).Value = '  +0.25%  '

$ws.Cells.Item(31, 5).Value = '  -0.24%  '

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '1.066'
$ws.Cells.Item(32, 5).Value = '  +6.44%  '

$ws.Cells.Item(33, 5).Value = '  +2.17%  '

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '1.419'
$ws.Cells.Item(34, 5).Value = '  +1.22%  '

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '3.542'
$ws.Cells.Item(35, 5).Value = '  -1.71%  '

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '5.378'
$ws.Cells.Item(36, 5).Value = '  +2.05%  '

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '0.06066'
$ws.Cells.Item(37, 5).Value = '  +0.33%  '

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.02243'
$ws.Cells.Item(38, 5).Value = '  +1.33%  '

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '8.212'
$ws.Cells.Item(39, 5).Value = '  -0.51%  '

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '1.170'
$ws.Cells.Item(40, 5).Value = '  -0.35%  '

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.5849'
$ws.Cells.Item(41, 5).Value = '  +1.42%  '

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '2.492'
$ws.Cells.Item(42, 5).Value = '  +10.48%  '

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '0.1833'
$ws.Cells.Item(43, 5).Value = '  +0.99%  '

$ws.Cells.Item(44, 5).Value = '  +0.23%  '

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '0.07870'
$ws.Cells.Item(45, 5).Value = '  +5.46%  '

$ws.Cells.Item(46, 5).Value = '  +1.12%  '

$ws.Cells.Item(47, 2).Value = 'Decentraland'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '0.5513'
$ws.Cells.Item(47, 5).Value = '  +1.29%  '

$ws.Cells.Item(48, 2).Value = 'EnergySwap'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '12.02'
$ws.Cells.Item(48, 5).Value = '  +0.81%  '

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '1.919'
$ws.Cells.Item(49, 5).Value = '  +1.57%  '

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '112.91'
$ws.Cells.Item(50, 5).Value = '  +1.63%  '

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '44.30'
$ws.Cells.Item(51, 5).Value = '  -2.20%  '
